# The commit swaps the "Integral" theme (ppt/theme/theme1.xml) and the
# "Office Theme" (ppt/theme/theme2.xml) between the slide master and the
# notes master: theme1.xml ends up holding the stock "Office Theme" colours
# (it was "Integral"), and theme2.xml ends up holding the "Integral"
# colours (it was "Office Theme").
#
# The presentation's only slide master is wired to theme1.xml, so we push
# the "Office Theme" palette onto it through the Theme colour scheme that
# is reachable from a Slide (this keeps the <a:clrScheme> element alive,
# unlike going through Master.ColorScheme). theme2.xml only backs the
# notes master and is not reachable/mutable through the exposed object
# model, so this focuses on the part of the swap that the host can apply.

function Get-ComRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme (target for theme1.xml / the slide master).
$officeTheme = @{
    1  = Get-ComRGB 0x00 0x00 0x00   # dk1
    2  = Get-ComRGB 0xFF 0xFF 0xFF   # lt1
    3  = Get-ComRGB 0x44 0x54 0x6A   # dk2
    4  = Get-ComRGB 0xE7 0xE6 0xE6   # lt2
    5  = Get-ComRGB 0x5B 0x9B 0xD5   # accent1
    6  = Get-ComRGB 0xED 0x7D 0x31   # accent2
    7  = Get-ComRGB 0xA5 0xA5 0xA5   # accent3
    8  = Get-ComRGB 0xFF 0xC0 0x00   # accent4
    9  = Get-ComRGB 0x44 0x72 0xC4   # accent5
    10 = Get-ComRGB 0x70 0xAD 0x47   # accent6
    11 = Get-ComRGB 0x05 0x63 0xC1   # hlink
    12 = Get-ComRGB 0x95 0x4F 0x72   # folHlink
}

# Apply through the first slide; every slide shares the single slide
# master/theme, so this repaints theme1.xml's <a:clrScheme> colours.
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeTheme[$i]
}
